$wb = $excel.ActiveWorkbook

# --- Update values on "PairCorrInput" sheet ---
$ws1 = $wb.Worksheets.Item("PairCorrInput")
$ws1.Activate() | Out-Null

$ws1.Range("E2").Value = 535
$ws1.Range("F2").Value = 90
$ws1.Range("F3").Value = 180
$ws1.Range("F4").Value = 365
$ws1.Range("F5").Value = 525

# New empty-but-formatted cells further down the sheet (F21, F22), carrying
# the same short-date formatting already used at K11/K12.
$ws1.Range("K11").Copy() | Out-Null
$ws1.Range("F21:F22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Selection on PairCorrInput moves to D11
$ws1.Range("D11").Select() | Out-Null

# --- Switch to the Information sheet, update its selection, and leave it
#     as the active (visible) tab, matching the saved workbook state. ---
$ws2 = $wb.Worksheets.Item("Information")
$ws2.Activate() | Out-Null
$ws2.Range("F33").Select() | Out-Null
